$d = $word.ActiveDocument

# Apply edits from the end of the document towards the start so that
# earlier (lower) character offsets are not invalidated by the
# length-changing edits performed later in the script.

# 9. "Brasilia-DF, 05 de junho de 2023." -> "06 de junho" (5 -> 6)
$r9 = $d.Range(846, 847)
if ($r9.Text -ne "5") { throw "r9 mismatch: [$($r9.Text)]" }
$r9.Text = "6"

# 8. Prepend new sentence before "Estamos abertos..." and lower-case the E
#    "Estamos abertos..." -> "A proposta ... e estamos abertos..."
$r8e = $d.Range(730, 731)
if ($r8e.Text -ne "E") { throw "r8e mismatch: [$($r8e.Text)]" }
$r8e.Text = "e"
$r8ins = $d.Range(730, 730)
$r8ins.InsertBefore("A proposta é condicionada à análise da documentação do imóvel e dos proprietários e ")

# 6. "um milhão e cinquenta " -> "um milhão e sessenta e cinco "
$r6 = $d.Range(477, 499)
if ($r6.Text -ne "um milhão e cinquenta ") { throw "r6 mismatch: [$($r6.Text)]" }
$r6.Text = "um milhão e sessenta e cinco "

# 5. "1.050" -> "1.065"
$r5 = $d.Range(463, 468)
if ($r5.Text -ne "1.050") { throw "r5 mismatch: [$($r5.Text)]" }
$r5.Text = "1.065"

# 4. "vinte" -> "quinze"
$r4 = $d.Range(365, 370)
if ($r4.Text -ne "vinte") { throw "r4 mismatch: [$($r4.Text)]" }
$r4.Text = "quinze"

# 3. "20" -> "15"
$r3 = $d.Range(354, 356)
if ($r3.Text -ne "20") { throw "r3 mismatch: [$($r3.Text)]" }
$r3.Text = "15"

# 2. "setenta" -> "oitenta"
$r2 = $d.Range(298, 305)
if ($r2.Text -ne "setenta") { throw "r2 mismatch: [$($r2.Text)]" }
$r2.Text = "oitenta"

# 1. "7" -> "8" (1.070.000,00 -> 1.080.000,00)
$r1 = $d.Range(275, 276)
if ($r1.Text -ne "7") { throw "r1 mismatch: [$($r1.Text)]" }
$r1.Text = "8"

Write-Output $d.Content.Text
